$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '307.82'
$ws.Range("E2").Value = '-2.91%'

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '37.96'
$ws.Range("E3").Value = '-3.86%'

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.050'
$ws.Range("E4").Value = '-1.70%'

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07893'
$ws.Range("E5").Value = '-3.58%'

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '2.021'
$ws.Range("E6").Value = '2.90%'

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '4.361'
$ws.Range("E7").Value = '2.75%'

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '8.212'
$ws.Range("E8").Value = '-0.23%'

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '3.124'
$ws.Range("E9").Value = '-1.20%'

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9263'
$ws.Range("E10").Value = '-0.25%'

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1278'
$ws.Range("E11").Value = '-9.54%'

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1902'
$ws.Range("E12").Value = '-4.78%'

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08742'
$ws.Range("E13").Value = '-3.05%'

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03433'
$ws.Range("E14").Value = '-1.84%'

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09735'
$ws.Range("E15").Value = '-0.91%'

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001399'
$ws.Range("E16").Value = '-0.42%'

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005981'
$ws.Range("E17").Value = '2.72%'

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D18").Value = '0.007506'
$ws.Range("E18").Value = '1,777.00%'

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '3.575'
$ws.Range("E19").Value = '-2.07%'

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3438'
$ws.Range("E20").Value = '-0.82%'

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '0.1284'
$ws.Range("E21").Value = '-1.49%'

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").Value = '5.017'
$ws.Range("E22").Value = '3.63%'

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '0.2517'
$ws.Range("E23").Value = '3.63%'

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = '0.04336'
$ws.Range("E24").Value = '-0.82%'

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D25").Value = '0.001221'
$ws.Range("E25").Value = '-0.15%'

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("B26").Value = 'HotbitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D26").Value = '0.004622'
$ws.Range("E26").Value = '-3.40%'

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("B27").Value = 'NitroEx'
$ws.Range("C27").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D27").Value = '0.0003593'
$ws.Range("E27").Value = '176.67%'

# Row 28
$ws.Range("B28").Value = 'Spectre.aiUtilityToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'

# Row 29
$ws.Range("B29").Value = 'LegolasExchange'
$ws.Range("C29").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'

# Row 30
$ws.Range("B30").Value = 'BitZToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'

# Row 31
$ws.Range("B31").Value = 'Birake'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'

# Row 32
$ws.Range("B32").Value = 'NashExchange'
$ws.Range("C32").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'

# Row 33
$ws.Range("B33").Value = 'AAXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'

# Row 34
$ws.Range("B34").Value = 'CenX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'

# Row 35
$ws.Range("B35").Value = 'BNIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02266'
$ws.Range("E39").Value = '3.05%'

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05030'
$ws.Range("E40").Value = '-2.83%'

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007602'
$ws.Range("E41").Value = '0.22%'

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009897'
$ws.Range("E42").Value = '1.06%'

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1362'
$ws.Range("E43").Value = '-0.86%'

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002028'
$ws.Range("E44").Value = '-4.69%'

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008555'
$ws.Range("E45").Value = '-6.35%'

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006414'
$ws.Range("E46").Value = '0.63%'

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000753'
$ws.Range("E47").Value = '0.40%'

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003007'
$ws.Range("E48").Value = '8.76%'

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001205'
$ws.Range("E49").Value = '0.41%'

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002108'
$ws.Range("E50").Value = '0.40%'

# Row 51
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002007'
$ws.Range("E51").Value = '0.40%'
